$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 124, shifting existing rows 124:149 down to 125:150.
$ws.Rows.Item(124).Insert()

# Populate the newly inserted row 124 with the new weekly price record.
# (Columns A,B,C,E,F,G,H,I,N,O,Q,R mirror the record that used to sit in
# row 124 before the shift; D,J,K,L,M,P carry the new week's figures.)
$ws.Cells.Item(124, 1).Value = 11
$ws.Cells.Item(124, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(124, 3).Value = "Bíobío"
$ws.Cells.Item(124, 4).Value = 45244
$ws.Cells.Item(124, 5).Value = 8
$ws.Cells.Item(124, 6).Value = 100112037
$ws.Cells.Item(124, 7).Value = "Cebollín"
$ws.Cells.Item(124, 8).Value = "Sin especificar"
$ws.Cells.Item(124, 9).Value = "Primera"
$ws.Cells.Item(124, 10).Value = 80
$ws.Cells.Item(124, 11).Value = 3000
$ws.Cells.Item(124, 12).Value = 3000
$ws.Cells.Item(124, 13).Value = 3000
$ws.Cells.Item(124, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(124, 15).Value = "Región Metropolitana"
$ws.Cells.Item(124, 16).Value = 83
$ws.Cells.Item(124, 17).Value = 36
$ws.Cells.Item(124, 18).Value = "Hortaliza"
